$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1053.2
$ws.Range("I15").Value = 1053.2
$ws.Range("K15").Value = 3159.6
$ws.Range("M15").Value = -2990.6
$ws.Range("H33").Value = 1148.409
$ws.Range("I33").Value = 1223.6842
$ws.Range("K33").Value = 1223.6842
$ws.Range("M33").Value = -994.6841999999999
$ws.Range("H107").Value = 981.375
$ws.Range("I107").Value = 900.36365
$ws.Range("J107").Value = 1159.6
$ws.Range("K107").Value = 900.36365
$ws.Range("L107").Value = 1159.6
$ws.Range("M107").Value = 1019.63635
$ws.Range("N107").Value = -4999.6
$ws.Range("H136").Value = 51631.58
$ws.Range("J136").Value = 51631.58
$ws.Range("L136").Value = 51631.58
$ws.Range("N136").Value = -61831.58
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5795.4395
$ws.Range("I32").Value = 2407.9348
$ws.Range("J32").Value = 13586.7
$ws.Range("K32").Value = 2407.9348
$ws.Range("L32").Value = 13586.7
$ws.Range("M32").Value = -2120.9348
$ws.Range("N32").Value = -14160.7
$ws.Range("H45").Value = 2777.28
$ws.Range("I45").Value = 1944.381
$ws.Range("K45").Value = 1944.381
$ws.Range("M45").Value = -1567.381
$ws.Range("H110").Value = 4029.2727
$ws.Range("I110").Value = 2924.889
$ws.Range("K110").Value = 2924.889
$ws.Range("M110").Value = -879.8890000000001
$ws.Range("H132").Value = 6131.9165
$ws.Range("I132").Value = 2337.0588
$ws.Range("K132").Value = 7011.176399999999
$ws.Range("M132").Value = -4481.176399999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -826
$ws.Range("H62").Value = 79981.5
$ws.Range("J62").Value = 79981.5
$ws.Range("L62").Value = 79981.5
$ws.Range("N62").Value = -81353.5
$ws.Range("H65").Value = 79981.5
$ws.Range("J65").Value = 79981.5
$ws.Range("L65").Value = 239944.5
$ws.Range("N65").Value = -246808.5
$ws.Range("H99").Value = 2894.75
$ws.Range("I99").Value = 2894.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2894.75
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -1396.75
$ws.Range("H105").Value = 4464.143
$ws.Range("I105").Value = 4116.3335
$ws.Range("K105").Value = 4116.3335
$ws.Range("M105").Value = -2369.3335
$ws.Range("H107").Value = 4470.4165
$ws.Range("I107").Value = 5169.65
$ws.Range("J107").Value = 974.25
$ws.Range("K107").Value = 5169.65
$ws.Range("L107").Value = 974.25
$ws.Range("M107").Value = -3249.65
$ws.Range("N107").Value = -4814.25
$ws.Range("H134").Value = 1824.8036
$ws.Range("I134").Value = 1262.0513
$ws.Range("J134").Value = 3115.8235
$ws.Range("K134").Value = 3786.1539
$ws.Range("L134").Value = 9347.470499999999
$ws.Range("M134").Value = -1251.1539
$ws.Range("N134").Value = -14417.4705
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2162.7104
$ws.Range("J31").Value = 3740
$ws.Range("L31").Value = 3740
$ws.Range("N31").Value = -4330
$ws.Range("H34").Value = 2162.7104
$ws.Range("J34").Value = 3740
$ws.Range("L34").Value = 3740
$ws.Range("N34").Value = -4144
$ws.Range("H58").Value = 2404.923
$ws.Range("I58").Value = 2073.4
$ws.Range("J58").Value = 2612.125
$ws.Range("K58").Value = 2073.4
$ws.Range("L58").Value = 2612.125
$ws.Range("M58").Value = -1870.4
$ws.Range("N58").Value = -3018.125
$ws.Range("H100").Value = 59971.75
$ws.Range("J100").Value = 59971.75
$ws.Range("L100").Value = 59971.75
$ws.Range("N100").Value = -62135.75
$ws.Range("H107").Value = 1008.1818
$ws.Range("I107").Value = 532.2
$ws.Range("K107").Value = 532.2
$ws.Range("M107").Value = 1387.8
$ws.Range("H131").Value = 78274.664
$ws.Range("J131").Value = 75329.60000000001
$ws.Range("L131").Value = 75329.60000000001
$ws.Range("N131").Value = -85409.60000000001
$ws.Range("H136").Value = 2404.923
$ws.Range("I136").Value = 2073.4
$ws.Range("J136").Value = 2612.125
$ws.Range("K136").Value = 6220.200000000001
$ws.Range("L136").Value = 7836.375
$ws.Range("M136").Value = -3670.200000000001
$ws.Range("N136").Value = -12936.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13865571
$ws.Range("I4").Value = 15845256
$ws.Range("J4").Value = 2647355
$ws.Range("K4").Value = 47535768
$ws.Range("L4").Value = 7942065
$ws.Range("M4").Value = -47535656
$ws.Range("N4").Value = -7942289
$ws.Range("H119").Value = 2515.889
$ws.Range("I119").Value = 2523.8333
$ws.Range("J119").Value = 2500
$ws.Range("K119").Value = 7571.499899999999
$ws.Range("L119").Value = 7500
$ws.Range("M119").Value = -2733.499899999999
$ws.Range("N119").Value = -17176
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 78916.664
$ws.Range("I43").Value = 36754
$ws.Range("J43").Value = 99998
$ws.Range("K43").Value = 36754
$ws.Range("L43").Value = 99998
$ws.Range("M43").Value = -36603
$ws.Range("N43").Value = -100300
$ws.Range("H46").Value = 3438.2
$ws.Range("I46").Value = 2547.75
$ws.Range("J46").Value = 7000
$ws.Range("K46").Value = 2547.75
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = -2391.75
$ws.Range("N46").Value = -7312
$ws.Range("H57").Value = 67447.5
$ws.Range("I57").Value = 20652.5
$ws.Range("J57").Value = 90845
$ws.Range("K57").Value = 20652.5
$ws.Range("L57").Value = 90845
$ws.Range("M57").Value = -19832.5
$ws.Range("N57").Value = -92485
$ws.Range("H80").Value = 6098.524
$ws.Range("I80").Value = 7461.077
$ws.Range("J80").Value = 3884.375
$ws.Range("K80").Value = 7461.077
$ws.Range("L80").Value = 3884.375
$ws.Range("M80").Value = -6463.077
$ws.Range("N80").Value = -5880.375
$ws.Range("H83").Value = 6098.524
$ws.Range("I83").Value = 7461.077
$ws.Range("J83").Value = 3884.375
$ws.Range("K83").Value = 37305.385
$ws.Range("L83").Value = 19421.875
$ws.Range("M83").Value = -32313.385
$ws.Range("N83").Value = -29405.875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 28827
$ws.Range("I56").Value = 25233.223
$ws.Range("J56").Value = 44999
$ws.Range("K56").Value = 25233.223
$ws.Range("L56").Value = 44999
$ws.Range("M56").Value = -24542.223
$ws.Range("N56").Value = -46381
$ws.Range("H102").Value = 94235
$ws.Range("J102").Value = 94235
$ws.Range("L102").Value = 94235
$ws.Range("N102").Value = -100725
$ws.Range("H132").Value = 5743.7896
$ws.Range("I132").Value = 5196.5
$ws.Range("J132").Value = 5808.1763
$ws.Range("K132").Value = 15589.5
$ws.Range("L132").Value = 17424.5289
$ws.Range("M132").Value = -13059.5
$ws.Range("N132").Value = -22484.5289
$ws.Range("H141").Value = 82928.5
$ws.Range("I141").Value = 86999
$ws.Range("J141").Value = 81571.664
$ws.Range("K141").Value = 86999
$ws.Range("L141").Value = 81571.664
$ws.Range("M141").Value = -81819
$ws.Range("N141").Value = -91931.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 3182.8
$ws.Range("I55").Value = 1199
$ws.Range("K55").Value = 1199
$ws.Range("M55").Value = -922
$ws.Range("H100").Value = 3871.25
$ws.Range("I100").Value = 1493.4
$ws.Range("K100").Value = 2986.8
$ws.Range("M100").Value = -2445.8
$ws.Range("H102").Value = 91244
$ws.Range("J102").Value = 91244
$ws.Range("L102").Value = 91244
$ws.Range("N102").Value = -97734
$ws.Range("H107").Value = 686.0625
$ws.Range("I107").Value = 532.5
$ws.Range("J107").Value = 805.5
$ws.Range("K107").Value = 1597.5
$ws.Range("L107").Value = 2416.5
$ws.Range("M107").Value = 322.5
$ws.Range("N107").Value = -6256.5
$ws.Range("H122").Value = 3805.3333
$ws.Range("I122").Value = 3939.6667
$ws.Range("J122").Value = 2999.3333
$ws.Range("K122").Value = 11819.0001
$ws.Range("L122").Value = 8997.999899999999
$ws.Range("M122").Value = -9369.000100000001
$ws.Range("N122").Value = -13897.9999
$ws.Range("H126").Value = 13507.517
$ws.Range("I126").Value = 9515.733
$ws.Range("K126").Value = 28547.199
$ws.Range("M126").Value = -26077.199
$ws.Range("H132").Value = 2266
$ws.Range("I132").Value = 2114.6155
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 6343.8465
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -3813.8465
$ws.Range("N132").Value = -14810
